$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (row 8) below the existing table.
$ws.Range("A8").Value = "Dilara "
$ws.Range("B8").Value = "aksoy"

# C8 looks like a date ("08.08.2022") but must stay literal text, matching
# the other rows in this column. A leading apostrophe forces text entry;
# resetting the style back to Normal afterwards avoids leaving a stray
# quote-prefix format on the cell.
$ws.Range("C8").Value = "'08.08.2022"
$ws.Range("C8").Style = "Normal"

$ws.Range("D8").Value = "Desk 046"
